$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '21.502.23'
$ws.Range("E2").Value = '  -2.71%  '

$ws.Range("D3").Value = "'" + '1.529.79'
$ws.Range("E3").Value = '  -1.80%  '

$ws.Range("D4").Value = "'" + '1.003'
$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D6").Value = "'" + '288.16'
$ws.Range("E6").Value = '  -1.47%  '

$ws.Range("D7").Value = "'" + '0.3872'
$ws.Range("E7").Value = '  -2.94%  '

$ws.Range("D8").Value = "'" + '0.3165'
$ws.Range("E8").Value = '  -2.19%  '

$ws.Range("D9").Value = "'" + '42.67'
$ws.Range("E9").Value = '  -3.11%  '

$ws.Range("D10").Value = "'" + '0.07144'
$ws.Range("E10").Value = '  -2.45%  '

$ws.Range("D11").Value = "'" + '1.066'
$ws.Range("E11").Value = '  -1.72%  '

$ws.Range("D12").Value = "'" + '1.003'
$ws.Range("E12").Value = '  +0.21%  '

$ws.Range("D13").Value = "'" + '5.717'
$ws.Range("E13").Value = '  +0.12%  '

$ws.Range("D14").Value = "'" + '18.12'
$ws.Range("E14").Value = '  -4.43%  '

$ws.Range("D15").Value = "'" + '6.536'
$ws.Range("E15").Value = '  -1.84%  '

$ws.Range("D16").Value = "'" + '1.535.06'
$ws.Range("E16").Value = '  -1.23%  '

$ws.Range("D17").Value = "'" + '0.00001086'
$ws.Range("E17").Value = '  -4.85%  '

$ws.Range("D18").Value = "'" + '0.06608'
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("D19").Value = "'" + '83.37'
$ws.Range("E19").Value = '  -0.60%  '

$ws.Range("D21").Value = "'" + '6.089'
$ws.Range("E21").Value = '  -3.57%  '

$ws.Range("D22").Value = "'" + '15.36'
$ws.Range("E22").Value = '  -2.48%  '

$ws.Range("D23").Value = "'" + '10.77'
$ws.Range("E23").Value = '  -4.78%  '

$ws.Range("D24").Value = "'" + '2.371'
$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").Value = "'" + '21.498.80'
$ws.Range("E25").Value = '  -2.76%  '

$ws.Range("D26").Value = "'" + '2.368'
$ws.Range("E26").Value = '  -3.46%  '

$ws.Range("D27").Value = "'" + '149.05'
$ws.Range("E27").Value = '  +0.23%  '

$ws.Range("E28").Value = '  -2.02%  '

$ws.Range("D29").Value = "'" + '4.831'
$ws.Range("E29").Value = '  -0.71%  '

$ws.Range("D30").Value = "'" + '1.708.14'
$ws.Range("E30").Value = '  -1.21%  '

$ws.Range("D31").Value = "'" + '116.41'
$ws.Range("E31").Value = '  -2.27%  '

$ws.Range("D32").Value = "'" + '6.029'
$ws.Range("E32").Value = '  +4.89%  '

$ws.Range("D33").Value = "'" + '0.9497'
$ws.Range("E33").Value = '  -6.25%  '

$ws.Range("D34").Value = "'" + '0.07998'
$ws.Range("E34").Value = '  -4.37%  '

$ws.Range("D35").Value = "'" + '8.476'
$ws.Range("E35").Value = '  -6.88%  '

$ws.Range("D36").Value = "'" + '5.150'
$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("D37").Value = "'" + '1.485'
$ws.Range("E37").Value = '  -8.53%  '

$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = "'" + '11.29'
$ws.Range("E38").Value = '  +4.65%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E39").Value = '  -3.56%  '

$ws.Range("D40").Value = "'" + '0.05868'
$ws.Range("E40").Value = '  -4.44%  '

$ws.Range("D41").Value = "'" + '0.2017'
$ws.Range("E41").Value = '  -2.23%  '

$ws.Range("D42").Value = "'" + '1.178'
$ws.Range("E42").Value = '  -3.35%  '

$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").Value = "'" + '0.5738'
$ws.Range("E44").Value = '  -2.29%  '

$ws.Range("D45").Value = "'" + '13.16'
$ws.Range("E45").Value = '  +0.76%  '

$ws.Range("D46").Value = "'" + '3.711'
$ws.Range("E46").Value = '  -1.50%  '

$ws.Range("D47").Value = "'" + '0.5544'
$ws.Range("E47").Value = '  -1.32%  '

$ws.Range("D48").Value = "'" + '1.888'
$ws.Range("E48").Value = '  -1.29%  '

$ws.Range("D49").Value = "'" + '1.158'
$ws.Range("E49").Value = '  +1.45%  '

$ws.Range("D50").Value = "'" + '115.36'
$ws.Range("E50").Value = '  -3.02%  '

$ws.Range("D51").Value = "'" + '0.06671'
$ws.Range("E51").Value = '  -2.67%  '
